# Swap the full contents (columns B through AD) between two rows,
# leaving column A (the running index) untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rowPairs = @(@(22, 23), @(142, 143))

foreach ($pair in $rowPairs) {
    $rowA = $pair[0]
    $rowB = $pair[1]

    # Columns B (2) .. AD (30)
    for ($col = 2; $col -le 30; $col++) {
        $cellA = $ws.Cells.Item($rowA, $col)
        $cellB = $ws.Cells.Item($rowB, $col)

        $valA = $cellA.Value2
        $valB = $cellB.Value2

        $cellA.Value = $valB
        $cellB.Value = $valA
    }
}
